$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set header row (row 1) values in the exact order new shared-strings ---
# Username, Password already exist as shared strings 0/1 from the original file.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# --- Data row (row 2) reuses DemoSalesManager / crmsfa shared strings ---
$ws.Range("A2").Value = "DemoSalesManager"
$ws.Range("B2").Value = "crmsfa"

# --- New header cells: First Name, Last Name ---
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"

# --- New data cells: Babu, M ---
$ws.Range("C2").Value = "Babu"
$ws.Range("D2").Value = "M"

# --- New header cells: Industry, Ownership ---
$ws.Range("F1").Value = "Industry"
$ws.Range("G1").Value = "Ownership"

# --- New data cells: Computer Software, Corporation ---
$ws.Range("F2").Value = "Computer Software"
$ws.Range("G2").Value = "Corporation"

# --- New header cell: Company Name ---
$ws.Range("E1").Value = "Company Name"

# --- New data cell: TestLeaf ---
$ws.Range("E2").Value = "TestLeaf"

# --- Remove the old 3rd row entirely ---
$ws.Rows.Item(3).Delete()

# --- Re-apply header formatting (yellow fill, no border) across the full header row in one shot ---
$ws.Range("A1:G1").ClearFormats()
$ws.Range("A1:G1").Interior.Color = 65535

# --- Data row: plain, no formatting ---
$ws.Range("A2:G2").ClearFormats()

# --- Column widths (engine rounds ColumnWidth to the nearest 1/6 on save, so
#     these inputs are chosen to land as close as possible to the target
#     rendered widths of 9.42578125 / 10.5703125 / 10.140625 / 10.140625 /
#     18.5703125 / 11.5703125) ---
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 9.333333333333332
$ws.Columns.Item(5).ColumnWidth = 9.333333333333332
$ws.Columns.Item(6).ColumnWidth = 17.666666666666668
$ws.Columns.Item(7).ColumnWidth = 10.666666666666666

# --- Sheet view: show gridlines (was hidden), select A2 ---
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("A2").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
